$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.927.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4679"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.75%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3707"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07380"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8730"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.48"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.808.23"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.370"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.84"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07074"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.511"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008726"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.963.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.338"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.061.51"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.903"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.61"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.222"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.34%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.334"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08937"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7701"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.497"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.908"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.64%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05287"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.958"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.334"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5353"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.370"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1670"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.466"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4962"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.47"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.675"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.94"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06298"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.44%  "
